# Add files via upload
#
# Sheet1's "cost" column (D) was populated with a literal placeholder
# string "NA" for every vegetable row. Replace it with the numeric
# value 0 for all 28 data rows (rows 2-29), turning the column into a
# real numeric column instead of text placeholders.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("D2:D29").Value = 0

# Leave the selection where the author ended up after editing the column.
$ws.Range("D2:D29").Select()
